$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend formatting (yellow header fill) to the new header cells R1:AE1
# by copying the format from the existing header cell Q1 before writing values.
$ws.Range("Q1").Copy()
$ws.Range("R1:AE1").PasteSpecial(-4122)

# --- Header row (row 1) values, column by column (A -> AE) ---
$ws.Range("A1").Value  = "Analyte_inchikey <anl>"
$ws.Range("B1").Value  = "Analyte_molfile"
$ws.Range("C1").Value  = "Analyte concentration [mM] <anl conc>"
$ws.Range("D1").Value  = "Solvent A_inchikey <solv>"
$ws.Range("E1").Value  = "Solvent A_molfile"
$ws.Range("F1").Value  = "Solvent volume [ml] <solv vol>"
$ws.Range("G1").Value  = "Additives <additives>"
$ws.Range("H1").Value  = "Additives concentration [mM] <additives conc>"
$ws.Range("I1").Value  = "Absorption Max [nm] <absorption_max>"
$ws.Range("J1").Value  = "Absorption Intensity Max [nm] <absorption_intensity_max>"
$ws.Range("K1").Value  = "Emission Max [nm] <emission_max>"
$ws.Range("L1").Value  = "Emission Intensity Max [nm] <emission_intensity_max>"
$ws.Range("M1").Value  = "Absorption Max 2 [nm] <absorption_max_2>"
$ws.Range("N1").Value  = "Absorption Intensity Max 2 [nm] <absorption_intensity_max_2>"
$ws.Range("O1").Value  = "Emission Max 2 [nm] <emission_max_2>"
$ws.Range("P1").Value  = "Emission Intensity Max 2 [nm] <emission_intensity_max_2>"
$ws.Range("Q1").Value  = "Absorption Max 3 [nm] <absorption_max_3>"
$ws.Range("R1").Value  = "Absorption Intensity Max 3 [nm] <absorption_intensity_max_3>"
$ws.Range("S1").Value  = "Emission Max 3 [nm] <emission_max_3>"
$ws.Range("T1").Value  = "Emission Intensity Max 3 [nm] <emission_intensity_max_3>"
$ws.Range("U1").Value  = "Interception number <interception_number>"
$ws.Range("V1").Value  = "Interception [nm] <interception>"
$ws.Range("W1").Value  = "Additives absorption [nm] <additives_absorption>"
$ws.Range("X1").Value  = "Intensity Additives absorption [nm] <intensity_absorption>"
$ws.Range("Y1").Value  = "Gas <gas>"
$ws.Range("Z1").Value  = "TemperatureP [°C] <temp>"
$ws.Range("AA1").Value = "E0 [eV] <auto-generated-E0>"
$ws.Range("AB1").Value = "Condition <condition>"
$ws.Range("AC1").Value = "Details <details>"
$ws.Range("AD1").Value = "Included <include>"
$ws.Range("AE1").Value = "BasePageName <BasePageName>"

# --- Column widths ---
# Column A: widen from the old best-fit width to an explicit custom width
# (the old sheet had bestFit on column A; the new layout uses a fixed,
# narrower, explicitly-set width instead).
$ws.Columns("A").ColumnWidth = 20.8
# New column AA gets its own explicit custom width (best-fit sized).
$ws.Columns("AA").ColumnWidth = 25.8

# --- Scroll the sheet view so column C is the first visible column,
# keeping the existing A2:XFD2 selection untouched.
$excel.ActiveWindow.ScrollColumn = 3
